$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 524.7222
$ws.Range("I41").Value = 274.42856
$ws.Range("J41").Value = 684
$ws.Range("K41").Value = 274.42856
$ws.Range("L41").Value = 684
$ws.Range("M41").Value = 165.57144
$ws.Range("N41").Value = -1564
$ws.Range("H129").Value = 1690.9375
$ws.Range("I129").Value = 818.1667
$ws.Range("J129").Value = 1892.3462
$ws.Range("K129").Value = 2454.5001
$ws.Range("L129").Value = 5677.0386
$ws.Range("M129").Value = 2545.4999
$ws.Range("N129").Value = -15677.0386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 914.61536
$ws.Range("I2").Value = 859
$ws.Range("J2").Value = 1100
$ws.Range("K2").Value = 859
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = -746
$ws.Range("N2").Value = -1326
$ws.Range("H5").Value = 87.71429000000001
$ws.Range("I5").Value = 85.666664
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 85.666664
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 26.333336
$ws.Range("N5").Value = -324
$ws.Range("H45").Value = 2966.6667
$ws.Range("I45").Value = 6000
$ws.Range("J45").Value = 1450
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 1450
$ws.Range("M45").Value = -5623
$ws.Range("N45").Value = -2204
$ws.Range("H57").Value = 18500
$ws.Range("I57").Value = 18500
$ws.Range("K57").Value = 18500
$ws.Range("M57").Value = -18016
$ws.Range("H61").Value = 2004.3948
$ws.Range("I61").Value = 2001.4445
$ws.Range("J61").Value = 2011.6364
$ws.Range("K61").Value = 2001.4445
$ws.Range("L61").Value = 2011.6364
$ws.Range("M61").Value = -1789.4445
$ws.Range("N61").Value = -2435.6364
$ws.Range("H74").Value = 3577.074
$ws.Range("I74").Value = 3645.423
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 3645.423
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = -2771.423
$ws.Range("N74").Value = -3548
$ws.Range("H77").Value = 3577.074
$ws.Range("I77").Value = 3645.423
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 18227.115
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = -13859.115
$ws.Range("N77").Value = -17736
$ws.Range("H116").Value = 914.61536
$ws.Range("I116").Value = 859
$ws.Range("J116").Value = 1100
$ws.Range("K116").Value = 859
$ws.Range("L116").Value = 1100
$ws.Range("M116").Value = 1435
$ws.Range("N116").Value = -5688
$ws.Range("H126").Value = 7666.6665
$ws.Range("I126").Value = 7666.6665
$ws.Range("K126").Value = 22999.9995
$ws.Range("M126").Value = -20529.9995
$ws.Range("H132").Value = 1495239.8
$ws.Range("I132").Value = 4335525
$ws.Range("J132").Value = 7471.143
$ws.Range("K132").Value = 13006575
$ws.Range("L132").Value = 22413.429
$ws.Range("M132").Value = -13004045
$ws.Range("N132").Value = -27473.429
$ws.Range("H136").Value = 2004.3948
$ws.Range("I136").Value = 2001.4445
$ws.Range("J136").Value = 2011.6364
$ws.Range("K136").Value = 6004.333500000001
$ws.Range("L136").Value = 6034.9092
$ws.Range("M136").Value = -3454.333500000001
$ws.Range("N136").Value = -11134.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 914.61536
$ws.Range("I3").Value = 859
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 859
$ws.Range("L3").Value = 1100
$ws.Range("M3").Value = -745
$ws.Range("N3").Value = -1328
$ws.Range("H4").Value = 87.71429000000001
$ws.Range("I4").Value = 85.666664
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 85.666664
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 29.333336
$ws.Range("N4").Value = -330
$ws.Range("H22").Value = 1221.2858
$ws.Range("I22").Value = 1137.4445
$ws.Range("J22").Value = 1724.3334
$ws.Range("K22").Value = 1137.4445
$ws.Range("L22").Value = 1724.3334
$ws.Range("M22").Value = -964.4445000000001
$ws.Range("N22").Value = -2070.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 1000
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 1000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -489
$ws.Range("N60").ClearContents()
$ws.Range("H105").Value = 1127.3928
$ws.Range("I105").Value = 1136.7693
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 1136.7693
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = 610.2307000000001
$ws.Range("N105").Value = -4499.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 103
$ws.Range("I8").Value = 103
$ws.Range("K8").Value = 309
$ws.Range("M8").Value = -170
$ws.Range("H103").Value = 5150
$ws.Range("I103").Value = 300
$ws.Range("J103").Value = 10000
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = -21
$ws.Range("N103").Value = -31758
$ws.Range("H105").Value = 10488.75
$ws.Range("J105").Value = 11269.143
$ws.Range("L105").Value = 33807.429
$ws.Range("N105").Value = -39049.429
$ws.Range("H106").Value = 4844.8
$ws.Range("J106").Value = 4844.8
$ws.Range("L106").Value = 14534.4
$ws.Range("N106").Value = -16426.4
$ws.Range("H121").Value = 3369228.2
$ws.Range("I121").Value = 431.66666
$ws.Range("J121").Value = 4117849.8
$ws.Range("K121").Value = 1294.99998
$ws.Range("L121").Value = 12353549.4
$ws.Range("M121").Value = 15.00001999999995
$ws.Range("N121").Value = -12356169.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 709.4828
$ws.Range("I107").Value = 771.5909
$ws.Range("J107").Value = 514.2857
$ws.Range("K107").Value = 771.5909
$ws.Range("L107").Value = 514.2857
$ws.Range("M107").Value = 1148.4091
$ws.Range("N107").Value = -4354.2857
$ws.Range("H132").Value = 7023
$ws.Range("I132").Value = 7274.6665
$ws.Range("J132").Value = 6915.143
$ws.Range("K132").Value = 21823.9995
$ws.Range("L132").Value = 20745.429
$ws.Range("M132").Value = -19293.9995
$ws.Range("N132").Value = -25805.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 38639.2
$ws.Range("J133").Value = 38639.2
$ws.Range("L133").Value = 38639.2
$ws.Range("N133").Value = -43699.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1678.5245
$ws.Range("I132").Value = 1546.8529
$ws.Range("J132").Value = 1844.3334
$ws.Range("K132").Value = 4640.5587
$ws.Range("L132").Value = 5533.0002
$ws.Range("M132").Value = -2110.5587
$ws.Range("N132").Value = -10593.0002
